# "add account and type_list to Ready to deploy"
#
# Two independent changes to the "基础数据" (base-data) sheet:
#
# 1. The sheet-wide "更新时间" (update-time) column (L) is bumped from
#    2021-06-05 to 2021-06-18 for every data row (rows 2-201).
#
# 2. A new/"ready to deploy" account ("一泡三响", previously sitting one
#    row below) is promoted ahead of "何泓姗VivaHo..." in the ranking
#    table - i.e. rows 89 and 90 swap their account-specific columns
#    (name / id / avatar / certification / fan-count / live-count /
#    product-count) while the row number (A), category (E), sales (J)
#    and period (K) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基础数据")

# xlPasteFormats - used below to restore the original "s=35" cell style
# after a NumberFormat round-trip (needed so purely-numeric-looking text
# like "37" / "146" or date-looking text like "2021-06-18" is written
# back as shared-string TEXT instead of being auto-coerced to a number).
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Bump the update-time column for every row (data rows are 2..201).
# ---------------------------------------------------------------------
$lastRow = 201
$dateRange = $ws.Range("L2:L" + $lastRow)
$dateRange.NumberFormat = "@"
$dateRange.Value = "2021-06-18"
$ws.Range("K2:K" + $lastRow).Copy()
$dateRange.PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 2) Swap the account data between row 89 (No.88) and row 90 (No.89).
# ---------------------------------------------------------------------
$row89 = @{
    B = "一泡三响"
    C = "1p3x"
    D = "//p26.douyinpic.com/aweme/1080x1080/31261000848cfae5cbf0a.webp?from=2956013662"
    F = "亳州市谯城区盛世华衣服饰官方账号"
    G = "545.71w"
    H = "37"
    I = "70"
}
$row90 = @{
    B = "何泓姗VivaHo今晚6:30直播"
    C = "Viva0410"
    D = "//p11.douyinpic.com/img/tos-cn-avt-0015/f165c94f2cadca6f00708fd841477458~c5_1080x1080.heic?from=2956013662"
    F = "演员、摄影师"
    G = "189.34w"
    H = "3"
    I = "146"
}

# Text-ish columns (account name / id / avatar url / certification /
# fan-count) are already non-numeric strings so a plain .Value write
# keeps them as shared-string text without any extra coaxing.
foreach ($col in @("B", "C", "D", "F", "G")) {
    $ws.Range($col + "89").Value = $row89[$col]
    $ws.Range($col + "90").Value = $row90[$col]
}

# H/I (live-stream count / product count) are plain digit strings
# ("37", "70", "3", "146") which Excel would otherwise auto-convert to
# numbers. Force text via NumberFormat, then restore the original
# "s=35" style via a format-only paste from an untouched same-style
# cell in the same row.
foreach ($col in @("H", "I")) {
    $cell89 = $ws.Range($col + "89")
    $cell89.NumberFormat = "@"
    $cell89.Value = $row89[$col]

    $cell90 = $ws.Range($col + "90")
    $cell90.NumberFormat = "@"
    $cell90.Value = $row90[$col]
}

$ws.Range("A89").Copy()
$ws.Range("H89:I89").PasteSpecial($xlPasteFormats)
$ws.Range("A90").Copy()
$ws.Range("H90:I90").PasteSpecial($xlPasteFormats)
